# Update market-price / profit figures across the Leve profit sheets
# (data refreshed by the scheduled Universalis price-sync runner)
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 69.36364
$ws.Range("I33").Value = 48.5
$ws.Range("K33").Value = 48.5
$ws.Range("M33").Value = 180.5
$ws.Range("H43").Value = 1559.6666
$ws.Range("I43").Value = 1199.5
$ws.Range("J43").Value = 1662.5714
$ws.Range("K43").Value = 1199.5
$ws.Range("L43").Value = 1662.5714
$ws.Range("M43").Value = -1130.5
$ws.Range("N43").Value = -1800.5714
$ws.Range("H58").Value = 794.1
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300
$ws.Range("H62").Value = 2977
$ws.Range("I62").Value = 2984
$ws.Range("K62").Value = 2984
$ws.Range("M62").Value = -2360
$ws.Range("H65").Value = 2977
$ws.Range("I65").Value = 2984
$ws.Range("K65").Value = 14920
$ws.Range("M65").Value = -11800
$ws.Range("H70").Value = 44080
$ws.Range("I70").Value = 40450
$ws.Range("J70").Value = 46500
$ws.Range("K70").Value = 121350
$ws.Range("L70").Value = 139500
$ws.Range("M70").Value = -121080
$ws.Range("N70").Value = -140040
$ws.Range("H73").Value = 44080
$ws.Range("I73").Value = 40450
$ws.Range("J73").Value = 46500
$ws.Range("K73").Value = 121350
$ws.Range("L73").Value = 139500
$ws.Range("M73").Value = -120414
$ws.Range("N73").Value = -141372
$ws.Range("H87").Value = 49900
$ws.Range("J87").Value = 49900
$ws.Range("L87").Value = 49900
$ws.Range("N87").Value = -52396
$ws.Range("H90").Value = 49900
$ws.Range("J90").Value = 49900
$ws.Range("L90").Value = 149700
$ws.Range("N90").Value = -162180
$ws.Range("H98").Value = 2729.1428
$ws.Range("I98").Value = 2729.1428
$ws.Range("K98").Value = 2729.1428
$ws.Range("M98").Value = -1231.1428
$ws.Range("H100").Value = 2499.4
$ws.Range("J100").Value = 2999
$ws.Range("L100").Value = 2999
$ws.Range("N100").Value = -4081
$ws.Range("H101").Value = 1433.3334
$ws.Range("J101").Value = 2000
$ws.Range("L101").Value = 6000
$ws.Range("N101").Value = -9244
$ws.Range("H116").Value = 13609.667
$ws.Range("I116").Value = 26574.25
$ws.Range("J116").Value = 3238
$ws.Range("K116").Value = 26574.25
$ws.Range("L116").Value = 3238
$ws.Range("M116").Value = -23132.25
$ws.Range("N116").Value = -10122
$ws.Range("H121").Value = 1498
$ws.Range("J121").Value = 1498
$ws.Range("L121").Value = 4494
$ws.Range("N121").Value = -7988
$ws.Range("H122").Value = 2729.1428
$ws.Range("I122").Value = 2729.1428
$ws.Range("K122").Value = 8187.428400000001
$ws.Range("M122").Value = -5737.428400000001
$ws.Range("H127").Value = 3454
$ws.Range("I127").Value = 3454
$ws.Range("K127").Value = 10362
$ws.Range("M127").Value = -5402

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 676.5
$ws.Range("I2").Value = 683.3333
$ws.Range("J2").Value = 656
$ws.Range("K2").Value = 683.3333
$ws.Range("L2").Value = 656
$ws.Range("M2").Value = -570.3333
$ws.Range("N2").Value = -882
$ws.Range("H32").Value = 3668.6
$ws.Range("I32").Value = 2211.4092
$ws.Range("K32").Value = 2211.4092
$ws.Range("M32").Value = -1924.4092
$ws.Range("H45").Value = 1722.5294
$ws.Range("I45").Value = 974
$ws.Range("K45").Value = 974
$ws.Range("M45").Value = -597
$ws.Range("H63").Value = 9999
$ws.Range("I63").Value = 9999
$ws.Range("K63").Value = 9999
$ws.Range("M63").Value = -9313
$ws.Range("H66").Value = 9999
$ws.Range("I66").Value = 9999
$ws.Range("K66").Value = 49995
$ws.Range("M66").Value = -46563
$ws.Range("H102").Value = 597.2
$ws.Range("I102").Value = 597.2
$ws.Range("K102").Value = 597.2
$ws.Range("M102").Value = 1024.8
$ws.Range("H116").Value = 676.5
$ws.Range("I116").Value = 683.3333
$ws.Range("J116").Value = 656
$ws.Range("K116").Value = 683.3333
$ws.Range("L116").Value = 656
$ws.Range("M116").Value = 1610.6667
$ws.Range("N116").Value = -5244
$ws.Range("H122").Value = 1602.4688
$ws.Range("I122").Value = 1570.5186
$ws.Range("K122").Value = 4711.5558
$ws.Range("M122").Value = -2261.5558

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 676.5
$ws.Range("I3").Value = 683.3333
$ws.Range("J3").Value = 656
$ws.Range("K3").Value = 683.3333
$ws.Range("L3").Value = 656
$ws.Range("M3").Value = -569.3333
$ws.Range("N3").Value = -884
$ws.Range("H80").Value = 7224.857
$ws.Range("I80").Value = 101.5
$ws.Range("K80").Value = 101.5
$ws.Range("M80").Value = 896.5
$ws.Range("H83").Value = 7224.857
$ws.Range("I83").Value = 101.5
$ws.Range("K83").Value = 507.5
$ws.Range("M83").Value = 4484.5
$ws.Range("H94").Value = 1174.75
$ws.Range("I94").Value = 966.3333
$ws.Range("J94").Value = 1800
$ws.Range("K94").Value = 966.3333
$ws.Range("L94").Value = 1800
$ws.Range("M94").Value = -515.3333
$ws.Range("N94").Value = -2702

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 987.25
$ws.Range("I16").Value = 916.3333
$ws.Range("K16").Value = 916.3333
$ws.Range("M16").Value = -629.3333
$ws.Range("H22").Value = 1009.75
$ws.Range("I22").Value = 269.5
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 269.5
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = 80.5
$ws.Range("N22").Value = -2450
$ws.Range("H113").Value = 987.25
$ws.Range("I113").Value = 916.3333
$ws.Range("K113").Value = 916.3333
$ws.Range("M113").Value = 1253.6667

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1888.6666
$ws.Range("I3").Value = 1099.75
$ws.Range("J3").Value = 3466.5
$ws.Range("K3").Value = 3299.25
$ws.Range("L3").Value = 10399.5
$ws.Range("M3").Value = -3187.25
$ws.Range("N3").Value = -10623.5
$ws.Range("H103").Value = 1871.1428
$ws.Range("J103").Value = 3799.5
$ws.Range("L103").Value = 11398.5
$ws.Range("N103").Value = -13156.5
$ws.Range("H131").Value = 13434.381
$ws.Range("I131").Value = 685.6667
$ws.Range("J131").Value = 14776.351
$ws.Range("K131").Value = 2057.0001
$ws.Range("L131").Value = 44329.053
$ws.Range("M131").Value = 2982.9999
$ws.Range("N131").Value = -54409.053

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5249.5
$ws.Range("I70").Value = 5624.25
$ws.Range("K70").Value = 5624.25
$ws.Range("M70").Value = -5354.25
$ws.Range("H73").Value = 5249.5
$ws.Range("I73").Value = 5624.25
$ws.Range("K73").Value = 5624.25
$ws.Range("M73").Value = -4688.25
$ws.Range("H102").Value = 1781.0834
$ws.Range("I102").Value = 1627.0769
$ws.Range("K102").Value = 1627.0769
$ws.Range("M102").Value = -5.076900000000023
$ws.Range("H132").Value = 3343
$ws.Range("I132").Value = 2709.077
$ws.Range("J132").Value = 6090
$ws.Range("K132").Value = 8127.231000000001
$ws.Range("L132").Value = 18270
$ws.Range("M132").Value = -5597.231000000001
$ws.Range("N132").Value = -23330

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6027.643
$ws.Range("J7").Value = 8236
$ws.Range("L7").Value = 8236
$ws.Range("N7").Value = -8460
$ws.Range("H40").Value = 9047.950000000001
$ws.Range("J40").Value = 10699.363
$ws.Range("L40").Value = 10699.363
$ws.Range("N40").Value = -10971.363
$ws.Range("H126").Value = 6027.643
$ws.Range("J126").Value = 8236
$ws.Range("L126").Value = 24708
$ws.Range("N126").Value = -29648

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4111.8945
$ws.Range("I136").Value = 4213.5
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 12640.5
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -10090.5
$ws.Range("N136").Value = -17097
